# Atualização automática de GUAIBA.xlsx
$wb = $excel.ActiveWorkbook

# 1) Rename "Paineis DARQ" -> "PAINEIS DARQ"
$wsPaineis = $wb.Worksheets.Item("Paineis DARQ")
$wsPaineis.Name = "PAINEIS DARQ"

# 2) Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$wsRecolhimento = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$wsRecolhimento.Name = "RECOLHIMENTO X ELIMINAÇÃO"

# 3) Update "DGC" sheet contents (keep its existing header/data formatting)
$wsDgc = $wb.Worksheets.Item("DGC")
$wsDgc.Range("A1").Value = "COMARCA"
$wsDgc.Range("B1").Value = "TEMÁTICA"
$wsDgc.Range("C1").Value = "PROBLEMA"
$wsDgc.Range("A2").Value = "Guaíba"
$wsDgc.Range("B2").Value = "MOT-VIG"
$wsDgc.Range("C2").Value = "Validou posto que não existe na comarca"

# 4) Delete the now-unused "Desarquivamentos Pendentes" sheet
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true
